$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (column index for subject / trial identifiers)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) updated meanEMG legmaxROM values
$ws.Range("B2").Value = 95.695904314898485
$ws.Range("C2").Value = 93.749991007127107
$ws.Range("D2").Value = 93.769939477364545
$ws.Range("E2").Value = 94.821777757501053

# Row 3 (STR) updated meanEMG legmaxROM values
$ws.Range("B3").Value = 94.174636240268143
$ws.Range("C3").Value = 94.067146827857201
$ws.Range("D3").Value = 91.702079432880268
$ws.Range("E3").Value = 94.798422624694567

# Selection now only spans the updated columns (B1:E3) instead of the full A Y range
$ws.Range("B1:E3").Select() | Out-Null
